$d = $word.ActiveDocument

# The document contains three "<id>...</id>" markers, each originally split
# across three separate runs: "<id>" (Courier New), the bare id text (plain),
# and "</id>" (Courier New). Collapse each trio into a single run whose text
# is the full "<id>...</id>" string, keeping the Courier New formatting that
# the "<id>" / "</id>" runs already used.

$ids = @("p015r_1", "p015r_2", "p015r_3")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $old, 2) | Out-Null
}
